$p = $ppt.ActivePresentation

# Locate the slide / shape that holds the old YouTube channel link in the
# "Rectangle 1" placeholder (the subtitle placeholder on the "Video link" slide).
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.Name -eq "Rectangle 1" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "https://www.youtube.com/@vidhibhanderi1308") {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
    }
}

$s = $targetSlide
$oldShape = $targetShape

# Remember the shape ids that exist before we touch anything.
$existingIds = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $existingIds += $s.Shapes.Item($i).Id
}

# Duplicate the shape first so the replacement keeps all of the original
# run/paragraph formatting (font, color, hyperlink relationship, etc.)
# byte-for-byte; PowerPoint hands the duplicate the next free shape id.
$dupRange = $oldShape.Duplicate()
$newShape = $dupRange.Item(1)

# Remove the original shape. Because it is tied to the slide layout's
# subtitle placeholder, PowerPoint immediately re-instantiates an empty
# placeholder shape in its place instead of just vanishing - clean up that
# leftover too so the slide ends up with exactly the shapes we want.
$oldShape.Delete()
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -ne $newShape.Id -and -not ($existingIds -contains $sh.Id)) {
        $sh.Delete()
    }
}

# Re-acquire a fresh reference to the duplicated shape (indices may have
# shifted after the deletes above).
$rect = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq $newShape.Id) {
        $rect = $sh
    }
}

# Update the text to the new video link.
$rect.TextFrame.TextRange.Text = "https://youtu.be/CxNWj2K19zE?si=zGJ-7p1rcW7D74x1"

# Nudge the shape to its new position/size (points; 12700 EMU per point).
$rect.Left = 1725613 / 12700.0
$rect.Top = 2287588 / 12700.0
$rect.Width = 9144000 / 12700.0
$rect.Height = 29.125
